$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
# Every cell referencing that string must be updated together so the report
# reads as freshly generated for handoff.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps bumped forward to reflect the new handoff generation time ---
$wsOverview.Range("G2").Value = "2016-08-28 07:05:24"
$wsDeDe.Range("H2").Value = "2016-08-28 07:05:24"
$wsZhCn.Range("H2").Value = "2016-08-28 07:05:20"

# --- Narrower status/date columns to fit the new, shorter content ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
